$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# Overview sheet: row 3 is for 906f1077-...md (header=row1, 888ccb3b=row2, 906f1077=row3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = $status
$wsZh.Range("G3").Value = "2016-03-09 22:49:53"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = $status
$wsDe.Range("G3").Value = "2016-03-09 22:50:08"
